# Work report: add a new entry row on Sheet1 (row 20, column B) containing
# the new value, mirroring a user typing into the next empty cell below the
# existing report rows and pressing Enter (which leaves the selection one
# row further down, on B21).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("B20").Value = "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"
$ws.Range("B21").Select()
